# Reorders the elements of the Python-set-literal strings stored in
# column N ("ref_transfer_all") for the rows whose set contents print in
# a different (but equivalent) order after the upstream function was
# revised. Only the text of these 59 cells changes; the underlying sets
# of ids are identical before/after.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = '{''47'', ''50'', ''54''}'
$ws.Range("N3").Value = '{''47'', ''50'', ''54''}'
$ws.Range("N6").Value = '{''54'', ''49'', ''47'', ''44'', ''66'', ''50''}'
$ws.Range("N9").Value = '{''50'', ''54'', ''53''}'
$ws.Range("N11").Value = '{''4'', ''50'', ''51'', ''62''}'
$ws.Range("N12").Value = '{''50'', ''51'', ''52''}'
$ws.Range("N13").Value = '{''24'', ''15'', ''62'', ''246'', ''51'', ''50'', ''2''}'
$ws.Range("N15").Value = '{''50'', ''51'', ''2''}'
$ws.Range("N16").Value = '{''50'', ''51'', ''17'', ''1''}'
$ws.Range("N18").Value = '{''13'', ''19'', ''7'', ''51'', ''50''}'
$ws.Range("N19").Value = '{''50'', ''51'', ''21'', ''7''}'
$ws.Range("N20").Value = '{''22'', ''231'', ''19'', ''369'', ''15'', ''61'', ''51'', ''50'', ''36''}'
$ws.Range("N23").Value = '{''22'', ''231'', ''19'', ''369'', ''15'', ''61'', ''51'', ''50'', ''36''}'
$ws.Range("N24").Value = '{''50'', ''51'', ''21'', ''7''}'
$ws.Range("N25").Value = '{''13'', ''19'', ''7'', ''51'', ''50''}'
$ws.Range("N27").Value = '{''50'', ''51'', ''17'', ''1''}'
$ws.Range("N28").Value = '{''50'', ''51'', ''2''}'
$ws.Range("N30").Value = '{''24'', ''15'', ''62'', ''246'', ''51'', ''50'', ''2''}'
$ws.Range("N31").Value = '{''50'', ''51'', ''52''}'
$ws.Range("N32").Value = '{''4'', ''50'', ''51'', ''62''}'
$ws.Range("N34").Value = '{''50'', ''54'', ''53''}'
$ws.Range("N37").Value = '{''54'', ''49'', ''47'', ''44'', ''66'', ''50''}'
$ws.Range("N40").Value = '{''47'', ''50'', ''54''}'
$ws.Range("N41").Value = '{''47'', ''50'', ''54''}'
$ws.Range("N43").Value = '{''54'', ''51'', ''53''}'
$ws.Range("N44").Value = '{''54'', ''240'', ''245'', ''53'', ''62'', ''37'', ''51'', ''65'', ''12'', ''40''}'
$ws.Range("N45").Value = '{''22'', ''14'', ''54'', ''24'', ''53'', ''13'', ''17'', ''51'', ''4'', ''12'', ''26'', ''52'', ''2'', ''43''}'
$ws.Range("N46").Value = '{''54'', ''51'', ''53''}'
$ws.Range("N47").Value = '{''14'', ''54'', ''51'', ''53''}'
$ws.Range("N48").Value = '{''54'', ''53'', ''19'', ''7'', ''1'', ''246'', ''51''}'
$ws.Range("N49").Value = '{''54'', ''51'', ''53''}'
$ws.Range("N50").Value = '{''54'', ''51'', ''53''}'
$ws.Range("N51").Value = '{''54'', ''240'', ''245'', ''53'', ''62'', ''37'', ''51'', ''65'', ''12'', ''40''}'
$ws.Range("N52").Value = '{''54'', ''51'', ''53''}'
$ws.Range("N53").Value = '{''54'', ''53'', ''19'', ''7'', ''1'', ''246'', ''51''}'
$ws.Range("N54").Value = '{''14'', ''54'', ''51'', ''53''}'
$ws.Range("N55").Value = '{''54'', ''51'', ''53''}'
$ws.Range("N56").Value = '{''22'', ''14'', ''54'', ''24'', ''53'', ''13'', ''17'', ''51'', ''4'', ''12'', ''26'', ''52'', ''2'', ''43''}'
$ws.Range("N57").Value = '{''245'', ''38'', ''34'', ''37'', ''52'', ''36''}'
$ws.Range("N58").Value = '{''52'', ''35'', ''34''}'
$ws.Range("N59").Value = '{''22'', ''14'', ''54'', ''24'', ''53'', ''13'', ''17'', ''51'', ''4'', ''12'', ''26'', ''52'', ''2'', ''43''}'
$ws.Range("N61").Value = '{''24'', ''19'', ''7'', ''1'', ''52''}'
$ws.Range("N64").Value = '{''50'', ''51'', ''52''}'
$ws.Range("N67").Value = '{''24'', ''19'', ''7'', ''1'', ''52''}'
$ws.Range("N69").Value = '{''22'', ''14'', ''54'', ''24'', ''53'', ''13'', ''17'', ''51'', ''4'', ''12'', ''26'', ''52'', ''2'', ''43''}'
$ws.Range("N70").Value = '{''52'', ''35'', ''34''}'
$ws.Range("N71").Value = '{''245'', ''38'', ''34'', ''37'', ''52'', ''36''}'
$ws.Range("N72").Value = '{''54'', ''51'', ''53''}'
$ws.Range("N73").Value = '{''50'', ''54'', ''53''}'
$ws.Range("N75").Value = '{''53'', ''44''}'
$ws.Range("N77").Value = '{''66'', ''41'', ''53''}'
$ws.Range("N78").Value = '{''47'', ''49'', ''41'', ''53''}'
$ws.Range("N79").Value = '{''47'', ''49'', ''53''}'
$ws.Range("N80").Value = '{''47'', ''49'', ''53''}'
$ws.Range("N81").Value = '{''47'', ''49'', ''41'', ''53''}'
$ws.Range("N82").Value = '{''66'', ''41'', ''53''}'
$ws.Range("N84").Value = '{''53'', ''44''}'
$ws.Range("N86").Value = '{''50'', ''54'', ''53''}'
$ws.Range("N87").Value = '{''54'', ''51'', ''53''}'
